$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to keep Text format so the literal string (with trailing zeros,
# leading zeros, etc.) is preserved exactly as in the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "35.558.99"
$ws.Range("D3").Value = "1.912.63"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "246.65"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  +5.61%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "42.08"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +5.38%  "
$ws.Range("D10").Value = "49.74"
$ws.Range("E10").Value = "  +6.55%  "
$ws.Range("D11").Value = "0.0721"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "2.190.55"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "12.32"
$ws.Range("E14").Value = "  +7.89%  "
$ws.Range("D15").Value = "0.701"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.913.81"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "4.91"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "35.580.78"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "72.50"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "0.0₃0824"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").Value = "245.77"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "12.68"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").Value = "4.85"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  +21.66%  "
$ws.Range("D27").Value = "171.26"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "8.43"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("D29").Value = "18.45"
$ws.Range("E29").Value = "  +4.60%  "
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "4.18"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("D32").Value = "0.0572"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "4.18"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("D35").Value = "0.922"
$ws.Range("E35").Value = "  +17.41%  "
$ws.Range("E36").Value = "  +5.03%  "
$ws.Range("D37").Value = "2.06"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").Value = "0.0212"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("D40").Value = "1.11"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("D41").Value = "0.0637"
$ws.Range("E41").Value = "  +14.80%  "
$ws.Range("D42").Value = "91.57"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "15.85"
$ws.Range("E43").Value = "  +7.51%  "
$ws.Range("D44").Value = "1.356.15"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +2.92%  "
$ws.Range("D46").Value = "47.59"
$ws.Range("E46").Value = "  +38.07%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "6.57"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "2.101.40"
$ws.Range("E51").Value = "  +3.35%  "
